$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new property (id 11) ---
$ws.Range("A12").Value = "11"
$ws.Range("B12").Value = "https://lh3.googleusercontent.com/pw/AP1GczMyb2ke1xP-GVhgkCSZ0U_aCMEt3elZtBxLS-aDohCKZ4MHpsGX_94vsLsBXVmFRPzLVQu5bQLm6pnJcIPFcwggz-8gefmjggC7SLZXkHqc9OMWi097xv8tpLsMlHhG78U3DH3TxHmuQc7ITz52i1jv=w1366-h616-s-no-gm?authuser=0"
$ws.Range("C12").Value = "lotes en la peña castilla la nueva "
$ws.Range("D12").Value = "lotes en barrio la peña, castilla la nueva. "
$ws.Range("E12").Value = '$37.000.000'
$ws.Range("F12").Value = "Castilla la Nueva"
$ws.Range("G12").Value = "3112697159"
$ws.Range("H12").Value = "0"
$ws.Range("I12").Value = "0"
$ws.Range("J12").Value = "https://lh3.googleusercontent.com/pw/AP1GczMyb2ke1xP-GVhgkCSZ0U_aCMEt3elZtBxLS-aDohCKZ4MHpsGX_94vsLsBXVmFRPzLVQu5bQLm6pnJcIPFcwggz-8gefmjggC7SLZXkHqc9OMWi097xv8tpLsMlHhG78U3DH3TxHmuQc7ITz52i1jv=w1366-h616-s-no-gm?authuser=0"

# --- Row 13: new property (id 12) ---
$ws.Range("A13").Value = "12"
$ws.Range("B13").Value = "https://lh3.googleusercontent.com/pw/AP1GczMtVnPP9S1H9RxqsK1GCYhJmVrfWHfXsE7xHstsx8750426Ansv6iRbwCLUIKiK-aQub2iIGO8k0cj4DFDZg36zzjCw_0toAO0hovy1dJUNURmVCuKm2MuTrUq6drMKaFdtuPmWY6NiPqbQZOu4J6o_=w1366-h616-s-no-gm?authuser=0"
$ws.Range("C13").Value = "casa Unifamiliar de dos pisos "
$ws.Range("D13").Value = "casa de 2 pisos Unifamiliar ubicada en castilla la nueva, barrio la peña. "
$ws.Range("E13").Value = '$350.000.000'
$ws.Range("F13").Value = "Castilla la Nueva"
$ws.Range("G13").Value = "3203441513"
$ws.Range("H13").Value = "4"
$ws.Range("I13").Value = "6"
$ws.Range("J13").Value = "https://lh3.googleusercontent.com/pw/AP1GczMtVnPP9S1H9RxqsK1GCYhJmVrfWHfXsE7xHstsx8750426Ansv6iRbwCLUIKiK-aQub2iIGO8k0cj4DFDZg36zzjCw_0toAO0hovy1dJUNURmVCuKm2MuTrUq6drMKaFdtuPmWY6NiPqbQZOu4J6o_=w1366-h616-s-no-gm?authuser=0"
$ws.Range("K13").Value = "https://lh3.googleusercontent.com/pw/AP1GczPSPCQQ7D1oRpABVoSf_K0yMHseMxJTT1nCMlQiPR4QzB5dmG24KuhZC4el1DHuCfbUcSxt4KD45xtoqfHwQTyTpj9lJHcI8Spk54pomSOfPSvdmff7ahXobTOnC14bBMUN0QMXoM_IwY1yd6sWiDuI=w289-h641-s-no-gm?authuser=0"

# --- Selection / view state: the new last row's contact photo cell ---
$ws.Range("K13").Select()
